$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the old scratch "wus" label in I6
# ---------------------------------------------------------------------
$ws.Range("I6").ClearContents()

# ---------------------------------------------------------------------
# 2) Block 1 (rows 8-13): replace the old I:K "1 test/2 test/3 test"
#    columns with the new J:M "wu_1..wu_4 runtime" columns, and add
#    AVERAGE formulas under C:E instead of the old K13 formula.
# ---------------------------------------------------------------------
$ws.Range("I8:K13").ClearContents()

$ws.Range("J8").Value2 = "wu_1 runtime"
$ws.Range("K8").Value2 = "wu_2 runtime"
$ws.Range("L8").Value2 = "wu_3 runtime"
$ws.Range("M8").Value2 = "wu_4 runtime"
$ws.Range("J8:M8").Style = "Heading 1"
$ws.Range("J8:M8").WrapText = $true
$ws.Range("J8:M8").Borders.Item(9).LineStyle = 0

$ws.Range("J9").Value2 = 81.900000000000006
$ws.Range("K9").Value2 = 81.900000000000006
$ws.Range("L9").Value2 = 81.900000000000006
$ws.Range("M9").Value2 = 81.900000000000006

$ws.Range("J10").Value2 = 87.52
$ws.Range("K10").Value2 = 87.52
$ws.Range("L10").Value2 = 87.53
$ws.Range("M10").Value2 = 87.53

$ws.Range("J11").Value2 = 82.13
$ws.Range("K11").Value2 = 82.13
$ws.Range("L11").Value2 = 87.6
$ws.Range("M11").Value2 = 87.6

$ws.Range("J9:M11").Font.Color = 0

$ws.Range("C13").Formula = "=AVERAGE(C9:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D9:D11)"
$ws.Range("E13").Formula = "=AVERAGE(E9:E11)"

# ---------------------------------------------------------------------
# 3) Block 2 (rows 20-25): add J:M header + data columns, add two new
#    rows of raw data (22, 23) and replace the old I/J AVERAGE row
#    with new C/D/E AVERAGE formulas in row 24. Old row 25 becomes
#    blank again.
# ---------------------------------------------------------------------
$ws.Range("I21:J25").ClearContents()

$ws.Range("J20").Value2 = "wu_1 runtime"
$ws.Range("K20").Value2 = "wu_2 runtime"
$ws.Range("L20").Value2 = "wu_3 runtime"
$ws.Range("M20").Value2 = "wu_4 runtime"
$ws.Range("J20:M20").Style = "Heading 1"
$ws.Range("J20:M20").WrapText = $true
$ws.Range("J20:M20").Borders.Item(9).LineStyle = 0

$ws.Range("J21").Value2 = 120.11
$ws.Range("K21").Value2 = 120.11
$ws.Range("L21").Value2 = 135.13
$ws.Range("M21").Value2 = 135.13

$ws.Range("C22").Value2 = 135.9
$ws.Range("D22").Value2 = 100.09
$ws.Range("E22").Value2 = 241.08
$ws.Range("J22").Value2 = 100.09
$ws.Range("K22").Value2 = 100.09
$ws.Range("L22").Value2 = 100.09
$ws.Range("M22").Value2 = 100.09

$ws.Range("C23").Value2 = 124.88
$ws.Range("D23").Value2 = 100.113
$ws.Range("E23").Value2 = 230
$ws.Range("J23").Value2 = 100.113
$ws.Range("K23").Value2 = 100.113
$ws.Range("L23").Value2 = 100.113
$ws.Range("M23").Value2 = 100.113

$ws.Range("C24").Formula = "=AVERAGE(C21:C23)"
$ws.Range("D24").Formula = "=AVERAGE(D21:D23)"
$ws.Range("E24").Formula = "=AVERAGE(E21:E23)"

# ---------------------------------------------------------------------
# 4) Block 4 (row 38 onward): add G/H "start time+qtime" columns and
#    J:N "wu runtime"/AVG columns, plus four new data rows 39-42.
# ---------------------------------------------------------------------
$ws.Range("G38").Value2 = "troy+diane start time+qtime"
$ws.Range("H38").Value2 = "troy+BJ start time+qtime"
$ws.Range("G38:H38").Style = "Heading 1"
$ws.Range("G38:H38").WrapText = $true

$ws.Range("J38").Value2 = "wu_1 runtime"
$ws.Range("K38").Value2 = "wu_2 runtime"
$ws.Range("L38").Value2 = "wu_3 runtime"
$ws.Range("M38").Value2 = "wu_4 runtime"
$ws.Range("N38").Value2 = "AVG"
$ws.Range("J38:N38").Style = "Heading 1"
$ws.Range("J38:N38").WrapText = $true
$ws.Range("J38:N38").Borders.Item(9).LineStyle = 0

$ws.Range("B39").Value2 = "1 node BJ + 1 node Diane"

$ws.Range("D39").Value2 = 90.84
$ws.Range("E39").Value2 = 241.56
$ws.Range("G39").Value2 = 136
$ws.Range("H39").Value2 = 83.7
$ws.Range("J39").Value2 = 100.07
$ws.Range("K39").Value2 = 100.07
$ws.Range("L39").Value2 = 78.989999999999995
$ws.Range("M39").Value2 = 84.26
$ws.Range("N39").Formula = "=AVERAGE(J39:M39)"

$ws.Range("D40").Value2 = 90.83
$ws.Range("E40").Value2 = 299.33999999999997
$ws.Range("G40").Value2 = 122.8
$ws.Range("H40").Value2 = 209.95
$ws.Range("J40").Value2 = 100.05
$ws.Range("K40").Value2 = 100.05
$ws.Range("L40").Value2 = 84.25
$ws.Range("M40").Value2 = 78.98
$ws.Range("N40").Formula = "=AVERAGE(J40:M40)"

$ws.Range("D41").Value2 = 92.16
$ws.Range("E41").Value2 = 215.9
$ws.Range("G41").Value2 = 110.65
$ws.Range("H41").Value2 = 84.46
$ws.Range("J41").Value2 = 100.059
$ws.Range("K41").Value2 = 100.059
$ws.Range("L41").Value2 = 84.265000000000001
$ws.Range("M41").Value2 = 84.265000000000001
$ws.Range("N41").Formula = "=AVERAGE(J41:M41)"

$ws.Range("D42").Value2 = 92.215000000000003
$ws.Range("E42").Value2 = 211.16
$ws.Range("G42").Value2 = 105.76
$ws.Range("H42").Value2 = 64.2
$ws.Range("J42").Value2 = 100.12
$ws.Range("K42").Value2 = 100.12
$ws.Range("L42").Value2 = 84.31
$ws.Range("M42").Value2 = 84.31
$ws.Range("N42").Formula = "=AVERAGE(J42:M42)"

# ---------------------------------------------------------------------
# 5) Sheet view: scroll so row 22 is near the top and select G44,
#    matching the author's final cursor position.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G44").Select()
